# Remove the post row for "「お金で買うことの出来ない物」" (row 725).
# Deleting the entire row shifts every subsequent row up by one, which
# matches the target diff (old row 726 becomes new row 725, ..., old
# row 775 becomes new row 774) and shrinks the used range/dimension
# from A1:C775 to A1:C774.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(725).Delete()
